$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values (C and D columns) for rows 2-7
$ws.Range("C2").Value = 0.03
$ws.Range("D2").Value = 0.01

$ws.Range("C3").Value = 0.19
$ws.Range("D3").Value = 0.02

$ws.Range("C4").Value = 0.97
$ws.Range("D4").Value = 0.03

$ws.Range("C5").Value = 2.88
$ws.Range("D5").Value = 0.05

$ws.Range("C6").Value = 5.99
$ws.Range("D6").Value = 0.07

$ws.Range("C7").Value = 10.64
$ws.Range("D7").Value = 0.11

# Clear rows 8 and 9 contents (A, C, D), and B formula too but keep formatting
$ws.Range("A8:D8").ClearContents()
$ws.Range("A9:D9").ClearContents()

# Re-apply shared formula for B2:B9 range
$ws.Range("B2:B9").Formula = "=A2*A2"
$ws.Range("B8:B9").ClearContents()

$ws.Range("H6").Select()
